$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the current "Tipo" column (E) values before shifting, since D/E will be
# overwritten with the new Latitud/Longitud columns and Tipo moves to column F.
$tipoHeader = $ws.Range("E1").Value2
$tipoValue  = $ws.Range("E2").Value2

# Rewrite header row: D/E become Latitud/Longitud, F becomes Tipo (moved from E)
$ws.Range("D1").Value = "Latitud"
$ws.Range("E1").Value = "Longitud"
$ws.Range("F1").Value = $tipoHeader

# Rewrite data row: D/E become numeric lat/long, F becomes the old Tipo value
$ws.Range("D2").Value = 15.26
$ws.Range("E2").Value = 26.56
$ws.Range("F2").Value = $tipoValue

$ws.Range("A1").Select()
$ws.Range("E2").Select()
